$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.943.74"

$ws.Range("D3").Value = "2.098.43"
$ws.Range("E3").Value = "  -1.19%  "

$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  -1.22%  "

$ws.Range("D5").Value = "'345.90"
$ws.Range("E5").Value = "  +2.25%  "

$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  -1.09%  "

$ws.Range("D7").Value = "'0.5143"
$ws.Range("E7").Value = "  -2.35%  "

$ws.Range("D8").Value = "'0.4419"
$ws.Range("E8").Value = "  -3.24%  "

$ws.Range("D9").Value = "'0.09364"
$ws.Range("E9").Value = "  +2.75%  "

$ws.Range("D10").Value = "'52.34"
$ws.Range("E10").Value = "  -4.86%  "

$ws.Range("D11").Value = "'1.168"
$ws.Range("E11").Value = "  -0.73%  "

$ws.Range("D12").Value = "'25.17"
$ws.Range("E12").Value = "  +2.74%  "

$ws.Range("D13").Value = "2.096.76"
$ws.Range("E13").Value = "  -1.04%  "

$ws.Range("D14").Value = "'6.732"
$ws.Range("E14").Value = "  -1.94%  "

$ws.Range("D15").Value = "'8.149"
$ws.Range("E15").Value = "  +0.07%  "

$ws.Range("D16").Value = "'99.48"
$ws.Range("E16").Value = "  +2.16%  "

$ws.Range("D17").Value = "'0.00001160"
$ws.Range("E17").Value = "  -1.38%  "

$ws.Range("E18").Value = "  -1.07%  "

$ws.Range("D19").Value = "'20.56"
$ws.Range("E19").Value = "  +5.14%  "

$ws.Range("D20").Value = "'0.06675"
$ws.Range("E20").Value = "  -0.51%  "

$ws.Range("D21").Value = "'1.001"
$ws.Range("E21").Value = "  -1.07%  "

$ws.Range("D22").Value = "'6.217"
$ws.Range("E22").Value = "  -2.24%  "

$ws.Range("D23").Value = "30.030.04"
$ws.Range("E23").Value = "  -2.44%  "

$ws.Range("D24").Value = "'12.59"
$ws.Range("E24").Value = "  -3.19%  "

$ws.Range("D25").Value = "'2.328"
$ws.Range("E25").Value = "  -1.58%  "

$ws.Range("D26").Value = "2.342.46"
$ws.Range("E26").Value = "  -1.10%  "

$ws.Range("D27").Value = "'21.99"
$ws.Range("E27").Value = "  -2.13%  "

$ws.Range("D28").Value = "'2.550"
$ws.Range("E28").Value = "  -0.32%  "

$ws.Range("D29").Value = "'162.52"
$ws.Range("E29").Value = "  -2.04%  "

$ws.Range("D30").Value = "'133.27"
$ws.Range("E30").Value = "  -1.62%  "

$ws.Range("D31").Value = "'1.165"
$ws.Range("E31").Value = "  -3.54%  "

$ws.Range("E32").Value = "  -1.87%  "

$ws.Range("D33").Value = "'1.637"
$ws.Range("E33").Value = "  -1.32%  "

$ws.Range("D34").Value = "'6.221"
$ws.Range("E34").Value = "  -2.69%  "

$ws.Range("D35").Value = "'3.947"
$ws.Range("E35").Value = "  -0.03%  "

$ws.Range("D36").Value = "'6.206"
$ws.Range("E36").Value = "  +5.10%  "

$ws.Range("E37").Value = "  -4.60%  "

$ws.Range("D38").Value = "'0.02559"
$ws.Range("E38").Value = "  -3.95%  "

$ws.Range("D39").Value = "'0.06771"
$ws.Range("E39").Value = "  -1.60%  "

$ws.Range("D40").Value = "'0.2277"
$ws.Range("E40").Value = "  -2.29%  "

$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.6922"
$ws.Range("E41").Value = "  -0.05%  "

$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").Value = "'12.48"
$ws.Range("E42").Value = "  -1.47%  "

$ws.Range("D43").Value = "'1.309"
$ws.Range("E43").Value = "  +3.68%  "

$ws.Range("D44").Value = "'0.6666"
$ws.Range("E44").Value = "  +2.54%  "

$ws.Range("D45").Value = "'14.24"
$ws.Range("E45").Value = "  -6.03%  "

$ws.Range("D46").Value = "'2.271"
$ws.Range("E46").Value = "  -1.92%  "

$ws.Range("D47").Value = "'3.627"
$ws.Range("E47").Value = "  -1.99%  "

$ws.Range("E48").Value = "  -5.86%  "

$ws.Range("D49").Value = "'1.220"
$ws.Range("E49").Value = "  -3.07%  "

$ws.Range("D50").Value = "'81.72"
$ws.Range("E50").Value = "  -2.42%  "

$ws.Range("D51").Value = "'0.07202"
$ws.Range("E51").Value = "  -1.53%  "
